# Auto update stocks_data.xlsx [2025-11-04 01:06:19]
#
# Adds a new column F holding the 2025/11/04 stock snapshot, mirroring the
# existing B:E layout (header date, header label, then alternating
# price/index rows with blank separator rows every third row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone column E's cell formatting (font, alignment, number format) onto
#    column F for every row so the new column matches the existing ones
#    before any values are written.
$ws.Range("E1:E115").Copy()
$ws.Range("F1:F115").PasteSpecial(-4122)

# 2) The full set of new column-F values, row by row (1-115). Blank
#    separator rows use $null so the cell stays empty (format-only), same
#    as the corresponding B/C/D/E cells in those rows.
$fValues = @(
    "2025/11/04",
    "上证",
    62.77,
    3976.52,
    $null,
    49.3,
    5594.41,
    $null,
    54.75,
    4653.4,
    $null,
    57.54,
    7333.6,
    $null,
    26.75,
    2715.84,
    $null,
    96.63,
    6851.97,
    $null,
    65.76000000000001,
    83978.49000000001,
    $null,
    85.7,
    19909.14,
    $null,
    83.77,
    39894.54,
    $null,
    58.13,
    5678.43,
    $null,
    11.57,
    33159.23,
    $null,
    29.11,
    3381.69,
    $null,
    47.4,
    3196.87,
    $null,
    18.72,
    7354.02,
    $null,
    31.86,
    8873.77,
    $null,
    13.4,
    13044.63,
    $null,
    24.18,
    12524.11,
    $null,
    21.59,
    9792.09,
    $null,
    26.86,
    16143.99,
    $null,
    32.17,
    17526.85,
    $null,
    20.84,
    10293.18,
    $null,
    15.19,
    9905.91,
    $null,
    20.88,
    3178.1,
    $null,
    43.62,
    5922.48,
    $null,
    29.08,
    9438.25,
    $null,
    13.29,
    2412.31,
    $null,
    56.14,
    2919.9,
    $null,
    58.79,
    3054.09,
    $null,
    52.38,
    3967.85,
    $null,
    47.08,
    2080.33,
    $null,
    28.02,
    14014.87,
    $null,
    86.36,
    9121.07,
    $null,
    57.02,
    12359.35,
    $null,
    6.08,
    2283.56,
    $null,
    26.06,
    879.7,
    $null,
    29.51,
    2779.5,
    $null,
    21.37,
    3998.7,
    $null,
    29.02,
    3372.75
)

# Write rows 2-115 directly (row 1 is a date-like string and needs special
# handling below to avoid being auto-converted to a date serial number).
for ($i = 1; $i -lt 115; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $fValues[$i]
}

# 3) Row 1 / column F holds "2025/11/04" which Excel would normally parse
#    as a date. Force Text format first so it is stored as a literal
#    string, matching the other header cells in row 1.
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = $fValues[0]

# Restore F1's visual formatting (it reverted to a generic Text style when
# NumberFormat was changed above) by re-copying column E's header format.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# 4) Give columns A-F an explicit width of 20 characters. Excel's
#    ColumnWidth property adds a small fixed padding before it is stored,
#    so 19.1666... (20 - 5/6) round-trips to a stored width of exactly 20.
$ws.Range("A1:F1").ColumnWidth = 19.166666666666668
